# Add a new "Run 50" column of data right before the existing "Mean"
# column. The old last column (AZ) keeps the "Run 50" header/values that
# used to describe the (now shifted) Mean numbers, and the Mean column
# itself moves one column to the right, into the new column BA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AZ1 becomes "Run 50", new BA1 becomes "Mean" (copy AZ1's
# formatting - bold font + border - onto the freshly used BA1 cell first)
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

# Data rows 2-14: AZ gets the new run values, BA gets the (shifted) mean
$newRunValue = 331.80192321
$newMeanValue = 265.70100603

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = $newRunValue   # column AZ = 52
    $ws.Cells.Item($r, 53).Value = $newMeanValue  # column BA = 53
}

Write-Output "done"
